# Updates the "relatorio_neomater_COMPLETO" worksheet so that it only
# reflects the current competencia (period) instead of cumulative data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 - PACOTE PRE-OPERATORIO PEDIATRICO CIRURGIA GERAL
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 8
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 4
$ws.Range("M3").Value = 4

# Row 10 - HERNIOPLASTIA INGUINAL (BILATERAL) - PEDIATRICO
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 4
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 2
$ws.Range("M10").Value = 0

# Row 11 - HERNIOPLASTIA UMBILICAL - PEDIATRICO
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 4
$ws.Range("M11").Value = 0

# Row 12 - ORQUIDOPEXIA BILATERAL - PEDIATRICO
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0

# Row 14 - CORRECAO DE HIPOSPADIA (1 TEMPO) - PEDIATRICO
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 1

# Row 16 - POSTECTOMIA - PEDIATRICO
$ws.Range("B16").Value = 4
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 4

# Row 17 - TOTAL
$ws.Range("B17").Value = 26
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 12
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 19
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 10
$ws.Range("M17").Value = 9
